$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (changed) date column C for rows 2-11 from 45243 (2023-11-13)
# to 45244 (2023-11-14), keeping the existing date formatting.
foreach ($row in 2..11) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45243) {
        $cell.Value = 45244
    }
}
